$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: "combien de dollars a chaque achat" goes from 1 to 2
$ws.Range("D2").Value = 2

# New buy entries logged in rows 8-11 (value, qty, date, time)
$ws.Range("A8").Value = 0.005505
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "20/11/2025"
$ws.Range("D8").Value = "19:02:10"

$ws.Range("A9").Value = 0.00539
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "21/11/2025"
$ws.Range("D9").Value = "01:01:07"

$ws.Range("A10").Value = 0.005185
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "21/11/2025"
$ws.Range("D10").Value = "07:01:30"

$ws.Range("A11").Value = 0.004827
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "21/11/2025"
$ws.Range("D11").Value = "14:02:09"

# Touch the sheet's very last row (mirrors the trailing formatted row stub
# left behind in the saved workbook) without disturbing any cell data.
$ws.Rows.Item(1048576).RowHeight = 12.75
